# Send mails to users every day
# Update stale sample e-mail addresses, fix a typo'd interest, and add the
# newest sign-up (Ivanna Parf) to the mailing list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refresh existing contacts -------------------------------------------------
$ws.Range("C2").Value = "08jccpq330@spymail.one"        # John Smith's e-mail
$ws.Range("C3").Value = "akkpeouddbvatl@dropmail.me"     # Marry Smith's e-mail
$ws.Range("D3").Value = "python"                          # Marry's interest
$ws.Range("C4").Value = "akkpeoomwzokhy@dropmail.me"     # Sim Kann's e-mail

# --- add the newest contact into the next free row -----------------------------
$ws.Range("A5").Value = "Ivanna"
$ws.Range("B5").Value = "Parf"
$ws.Range("C5").Value = "mainpy571@gmail.com"
$ws.Range("D5").Value = "yoga"
